$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new abbreviation rows after the last existing row (144)
$ws.Range("A145").Value = "tPA"
$ws.Range("B145").Value = "Tissue plasminogen activator"
$ws.Range("A146").Value = "PAF"
$ws.Range("B146").Value = "Platelet-activating factor"

# Update the view: scroll down and move the selection to the new data
$ws.Application.ActiveWindow.ScrollRow = 132
$ws.Range("A145").Select()
